# Applies updated win-probability matrix values for Saint Mary's (CA)_A
# team-specific matrix, based on games pulled March 7.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1928104575163399
$ws.Range("C2").Value = 0.565359477124183
$ws.Range("J2").Value = 0.02287581699346405
$ws.Range("P2").Value = 0.1405228758169935
$ws.Range("S2").Value = 0.07843137254901961
$ws.Range("B3").Value = 0.005714285714285714
$ws.Range("C3").Value = 0.01714285714285714
$ws.Range("J3").Value = 0.02285714285714286
$ws.Range("P3").Value = 0.76
$ws.Range("S3").Value = 0.1942857142857143
$ws.Range("J4").Value = 0.03174603174603174
$ws.Range("P4").Value = 0.7619047619047619
$ws.Range("S4").Value = 0.2063492063492063
$ws.Range("B6").Value = 0.07172995780590717
$ws.Range("D6").Value = 0.01265822784810127
$ws.Range("F6").Value = 0.07172995780590717
$ws.Range("J6").Value = 0.2067510548523207
$ws.Range("O6").Value = 0.004219409282700422
$ws.Range("Q6").Value = 0.1983122362869198
$ws.Range("R6").Value = 0.0759493670886076
$ws.Range("S6").Value = 0.3586497890295359
$ws.Range("B7").Value = 0.1165644171779141
$ws.Range("D7").Value = 0.01840490797546012
$ws.Range("F7").Value = 0.05521472392638037
$ws.Range("J7").Value = 0.1595092024539877
$ws.Range("O7").Value = 0.03067484662576687
$ws.Range("Q7").Value = 0.2269938650306748
$ws.Range("R7").Value = 0.03680981595092025
$ws.Range("S7").Value = 0.3558282208588957
$ws.Range("B8").Value = 0.09662921348314607
$ws.Range("D8").Value = 0.01797752808988764
$ws.Range("F8").Value = 0.07191011235955057
$ws.Range("J8").Value = 0.1348314606741573
$ws.Range("O8").Value = 0.01348314606741573
$ws.Range("Q8").Value = 0.2022471910112359
$ws.Range("R8").Value = 0.07191011235955057
$ws.Range("S8").Value = 0.3910112359550562
$ws.Range("B9").Value = 0.1085714285714286
$ws.Range("D9").Value = 0.05142857142857143
$ws.Range("E9").Value = 0.005714285714285714
$ws.Range("F9").Value = 0.09714285714285714
$ws.Range("J9").Value = 0.09714285714285714
$ws.Range("O9").Value = 0.02857142857142857
$ws.Range("Q9").Value = 0.1542857142857143
$ws.Range("R9").Value = 0.09714285714285714
$ws.Range("S9").Value = 0.36
$ws.Range("B10").Value = 0.1074870274277242
$ws.Range("D10").Value = 0.03187546330615271
$ws.Range("E10").Value = 0.002223869532987398
$ws.Range("F10").Value = 0.07116382505559674
$ws.Range("J10").Value = 0.1326908821349148
$ws.Range("O10").Value = 0.01779095626389918
$ws.Range("Q10").Value = 0.2164566345441067
$ws.Range("R10").Value = 0.08673091178650852
$ws.Range("S10").Value = 0.3335804299481097
$ws.Range("F11").Value = 0.003861003861003861
$ws.Range("G11").Value = 0.1467181467181467
$ws.Range("J11").Value = 0.1003861003861004
$ws.Range("K11").Value = 0.2200772200772201
$ws.Range("L11").Value = 0.5096525096525096
$ws.Range("S11").Value = 0.0193050193050193
$ws.Range("G12").Value = 0.7279411764705882
$ws.Range("J12").Value = 0.1911764705882353
$ws.Range("K12").Value = 0.02205882352941177
$ws.Range("L12").Value = 0.02205882352941177
$ws.Range("S12").Value = 0.03676470588235294
$ws.Range("G13").Value = 0.6122448979591837
$ws.Range("J13").Value = 0.3061224489795918
$ws.Range("S13").Value = 0.08163265306122448
$ws.Range("F15").Value = 0.01433691756272401
$ws.Range("H15").Value = 0.1362007168458781
$ws.Range("I15").Value = 0.06810035842293907
$ws.Range("J15").Value = 0.4050179211469534
$ws.Range("K15").Value = 0.05376344086021505
$ws.Range("M15").Value = 0.007168458781362007
$ws.Range("N15").Value = 0.003584229390681004
$ws.Range("S15").Value = 0.2007168458781362
$ws.Range("F16").Value = 0.009174311926605505
$ws.Range("H16").Value = 0.1697247706422018
$ws.Range("I16").Value = 0.08256880733944955
$ws.Range("J16").Value = 0.4036697247706422
$ws.Range("K16").Value = 0.1146788990825688
$ws.Range("M16").Value = 0.004587155963302753
$ws.Range("O16").Value = 0.0871559633027523
$ws.Range("S16").Value = 0.1284403669724771
$ws.Range("F17").Value = 0.006198347107438017
$ws.Range("H17").Value = 0.1735537190082645
$ws.Range("I17").Value = 0.07231404958677685
$ws.Range("J17").Value = 0.4669421487603306
$ws.Range("K17").Value = 0.08057851239669421
$ws.Range("M17").Value = 0.01033057851239669
$ws.Range("O17").Value = 0.06611570247933884
$ws.Range("S17").Value = 0.1239669421487603
$ws.Range("F18").Value = 0.01595744680851064
$ws.Range("H18").Value = 0.1436170212765958
$ws.Range("I18").Value = 0.07446808510638298
$ws.Range("J18").Value = 0.4680851063829787
$ws.Range("K18").Value = 0.0851063829787234
$ws.Range("M18").Value = 0.005319148936170213
$ws.Range("N18").Value = 0.005319148936170213
$ws.Range("O18").Value = 0.06382978723404255
$ws.Range("S18").Value = 0.1382978723404255
$ws.Range("F19").Value = 0.01897689768976898
$ws.Range("H19").Value = 0.2120462046204621
$ws.Range("I19").Value = 0.07508250825082509
$ws.Range("J19").Value = 0.363036303630363
$ws.Range("K19").Value = 0.08333333333333333
$ws.Range("M19").Value = 0.03465346534653466
$ws.Range("O19").Value = 0.09075907590759076
$ws.Range("S19").Value = 0.1221122112211221
